$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.000" or "312.84"
# are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Each entry: Row, Coin (B), Link (C), Price (D), Volume1h (E)
$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '27.563.99', '  -4.81%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.841.32', '  -4.01%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.000', '  -0.58%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '312.84', '  -3.88%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9994', '  -0.56%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4248', '  -7.30%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3627', '  -4.75%  '),
    @(9, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '43.65', '  -4.25%  '),
    @(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07209', '  -6.90%  '),
    @(11, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.8960', '  -8.17%  '),
    @(12, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '20.63', '  -8.15%  '),
    @(13, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.828.85', '  -6.97%  '),
    @(14, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.573', '  -5.51%  '),
    @(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.321', '  -6.56%  '),
    @(16, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06819', '  -2.36%  '),
    @(17, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '0.9996', '  -0.75%  '),
    @(18, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '77.36', '  -8.55%  '),
    @(19, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000008895', '  -6.09%  '),
    @(20, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9994', '  -0.51%  '),
    @(21, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '15.36', '  -7.83%  '),
    @(22, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '27.528.51', '  -4.97%  '),
    @(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.942', '  -7.42%  '),
    @(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '10.69', '  -3.38%  '),
    @(25, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.068.71', '  -4.55%  '),
    @(26, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.047', '  -0.58%  '),
    @(27, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '151.68', '  -3.95%  '),
    @(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '18.15', '  -4.79%  '),
    @(29, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '5.307', '  -5.36%  '),
    @(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '110.72', '  -5.82%  '),
    @(31, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.731', '  -5.69%  '),
    @(32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08881', '  -4.27%  '),
    @(33, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7752', '  -10.13%  '),
    @(34, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.468', '  -12.45%  '),
    @(35, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.856', '  -5.19%  '),
    @(36, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.082', '  -12.56%  '),
    @(37, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '0.9988', '  -0.65%  '),
    @(38, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05395', '  -5.07%  '),
    @(39, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.096', '  -4.74%  '),
    @(40, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.962', '  -3.78%  '),
    @(41, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01918', '  -6.00%  '),
    @(42, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5041', '  -8.19%  '),
    @(43, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.785', '  -8.97%  '),
    @(44, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1630', '  -7.13%  '),
    @(45, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.06614', '  -4.59%  '),
    @(46, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '8.210', '  -11.77%  '),
    @(47, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '105.70', '  -4.33%  '),
    @(48, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.4699', '  -9.12%  '),
    @(49, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '10.23', '  -8.30%  '),
    @(50, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9986', '  -0.63%  '),
    @(51, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.641', '  -6.76%  ')
)

foreach ($item in $data) {
    $row = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
}
